$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix minor coordinate errors on map for four farms
# (order matters for shared-string table insertion order)
$ws.Range("B47").Value = "47.447300, -122.459900"
$ws.Range("B43").Value = "47.2201236,-121.9790837"
$ws.Range("B33").Value = "47.4746452,-122.2849894"
$ws.Range("B34").Value = "47.4763239,-122.2807105"

# B47's new value is longer, so make it readable (left/middle aligned, wrapped)
$r47 = $ws.Range("B47")
$r47.Font.Size = 10
$r47.Font.Color = 0
$r47.Font.Name = "Inherit"
$r47.HorizontalAlignment = -4131
$r47.VerticalAlignment = -4108
$r47.WrapText = $true

# Scroll the view to where the edits were made
[void]$ws.Range("B34").Select()

$ws.PageSetup.Orientation = 1
